$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 917, shifting existing rows 917:963 down to 918:964
$ws.Rows.Item(917).Insert()

# Populate the newly inserted row 917 with the new weekly record
$ws.Cells.Item(917, 1).Value = 10
$ws.Cells.Item(917, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(917, 3).Value = "La Araucanía"
$ws.Cells.Item(917, 4).Value = 45267
$ws.Cells.Item(917, 5).Value = 9
$ws.Cells.Item(917, 6).Value = 100112043
$ws.Cells.Item(917, 7).Value = "Pepino ensalada"
$ws.Cells.Item(917, 8).Value = "Sin especificar"
$ws.Cells.Item(917, 9).Value = "Primera"
$ws.Cells.Item(917, 10).Value = 500
$ws.Cells.Item(917, 11).Value = 16000
$ws.Cells.Item(917, 12).Value = 19000
$ws.Cells.Item(917, 13).Value = 16600
$ws.Cells.Item(917, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(917, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(917, 16).Value = 332
$ws.Cells.Item(917, 17).Value = 50
$ws.Cells.Item(917, 18).Value = "Hortaliza"
